$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: merge "AirDensityDrive" + ": F1 Power & Top-
#    Speed Analysis" runs into a single run and drop the spell-check
#    proofErr markers that wrapped "AirDensityDrive".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Executive Summary of AirDensityDrive: F1 Power & Top-Speed Analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Executive Summary of AirDensityDrive: F1 Power & Top-Speed Analysis", 2)

# ------------------------------------------------------------------
# 2) "Top Speed Solver" bullet: merge the first four runs into one
#    ("Top Speed Solver: Numerical root-finding(fzero) to balance
#    Pdrag") and strip every proofErr marker in the paragraph, while
#    leaving the remaining runs (" ", "+", " ", "P", "rr", " ", "=",
#    " ", "Ptotal") exactly as they were, only without the proofErr
#    wrappers. Rebuilding the whole paragraph via InsertXML gives an
#    exact, predictable result (Find/Replace alone tends to merge
#    more runs than intended and can orphan proofErr markers).
# ------------------------------------------------------------------
$solverPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Top Speed Solver")) {
        $solverPara = $cand
    }
}

$solverXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Top Speed Solver: Numerical root-finding(fzero) to balance Pdrag</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>P</w:t></w:r><w:r><w:t>rr</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Ptotal</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$solverPara.Range.InsertXML($solverXml)

# ------------------------------------------------------------------
# 3) Add a new "Discrepancies" paragraph right after the "Plots"
#    paragraph.
# ------------------------------------------------------------------
$plotsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Trim() -eq "Plots") {
        $plotsPara = $cand
        $plotsIndex = $i
    }
}

$plotsPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($plotsIndex + 1)
$newPara.Range.Text = "Discrepancies"
